$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.102335785814724
$ws.Range("C2").Value = 6.135771264509594

$ws.Range("B3").Value = 6.090649959166668
$ws.Range("C3").Value = 6.169074292372931

$ws.Range("B4").Value = 6.137456155375406

$ws.Range("B5").Value = 5.930117065756425
$ws.Range("C5").Value = 5.744895021188859

$ws.Range("B6").Value = 5.983852791196446
$ws.Range("C6").Value = 6.090875258850894

$ws.Range("B7").Value = 6.307698047608095
$ws.Range("C7").Value = 5.775077483633254

$ws.Range("B8").Value = 5.465312693427742
$ws.Range("C8").Value = 6.484699172290198

$ws.Range("B9").Value = 6.240186078492074
$ws.Range("C9").Value = 6.175539651802586

$ws.Range("B10").Value = 6.271727777662242
$ws.Range("C10").Value = 5.299931936610919

$ws.Range("B11").Value = 5.65361300174162
$ws.Range("C11").Value = 6.399498130247556
